$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1725
$ws1.Range("F6").Value = 624
$ws1.Range("F7").Value = 1130
$ws1.Range("F8").Value = 1564
$ws1.Range("F9").Value = 168
$ws1.Range("F11").Value = 1483
$ws1.Range("F12").Value = 3112
$ws1.Range("F13").Value = 660
$ws1.Range("F14").Value = 1801
$ws1.Range("F15").Value = 1812
$ws1.Range("F16").Value = 864
$ws1.Range("F17").Value = 289
$ws1.Range("F18").Value = 9
$ws1.Range("F19").Value = 1492
$ws1.Range("F20").Value = 294
$ws1.Range("F21").Value = 73
$ws1.Range("F22").Value = 20
$ws1.Range("F23").Value = 1251
$ws1.Range("F25").Value = 475
$ws1.Range("F26").Value = 150
$ws1.Range("F27").Value = 4876
$ws1.Range("F28").Value = 5272
$ws1.Range("F29").Value = 762
$ws1.Range("F30").Value = 585
$ws1.Range("F31").Value = 1679
$ws1.Range("F33").Value = 186

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 62
$ws2.Range("F4").Value = 86
$ws2.Range("F9").Value = 103

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 45

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 45
$ws4.Range("F6").Value = 62
$ws4.Range("F7").Value = 86
$ws4.Range("F9").Value = 1725
$ws4.Range("F11").Value = 624
$ws4.Range("F12").Value = 1130
$ws4.Range("F13").Value = 1564
$ws4.Range("F14").Value = 168
$ws4.Range("F15").Value = 168
$ws4.Range("F18").Value = 1483
$ws4.Range("F19").Value = 3112
$ws4.Range("F20").Value = 660
$ws4.Range("F21").Value = 1801
$ws4.Range("F22").Value = 1812
$ws4.Range("F23").Value = 864
$ws4.Range("F24").Value = 289
$ws4.Range("F25").Value = 9
$ws4.Range("F26").Value = 1492
$ws4.Range("F27").Value = 294
$ws4.Range("F28").Value = 73
$ws4.Range("F30").Value = 20
$ws4.Range("F32").Value = 1251
$ws4.Range("F34").Value = 475
$ws4.Range("F35").Value = 150
$ws4.Range("F36").Value = 4876
$ws4.Range("F37").Value = 5272
$ws4.Range("F38").Value = 762
$ws4.Range("F39").Value = 585
$ws4.Range("F40").Value = 1679
$ws4.Range("F41").Value = 103
$ws4.Range("F44").Value = 186
